$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Reactivos")

# The reagent "Activo" field used to be bound to a description property;
# repoint the template placeholder to the simpler {{Reactivo.Activo}} field.
$ws.Range("B11").Value = "{{Reactivo.Activo}}"

# Restore the active selection to G1 (matches the latest save state).
[void]$ws.Range("G1").Select()
